$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(78, 48, 2024),
    @(79, 45, 2024),
    @(80, 26, 2024),
    @(81, 45, 2024),
    @(82, 44, 2024),
    @(83, 43, 2024),
    @(84, 41, 2024),
    @(85, 41, 2024),
    @(86, 40, 2024),
    @(87, 39, 2024),
    @(88, 35, 2024),
    @(89, 24, 2024),
    @(90, 36, 2024),
    @(91, 36, 2024),
    @(92, 36, 2024),
    @(93, 36, 2024),
    @(94, 30, 2024),
    @(95, 29, 2024),
    @(96, 29, 2024),
    @(97, 29, 2024),
    @(98, 28, 2024),
    @(99, 28, 2024),
    @(100, 28, 2024),
    @(101, 22, 2024),
    @(102, 23, 2024),
    @(103, 23, 2024),
    @(104, 23, 2024),
    @(105, 20, 2024),
    @(106, 16, 2024),
    @(107, 11, 2024),
    @(108, 8, 2024),
    @(109, 7, 2024),
    @(110, 6, 2024),
    @(111, 6, 2024),
    @(112, 6, 2024),
    @(113, 6, 2024),
    @(114, 7, 2024),
    @(115, 7, 2024),
    @(116, 6, 2024),
    @(117, 6, 2024),
    @(118, 2, 2024),
    @(119, 7, 2024),
    @(120, 4, 2024),
    @(121, 7, 2024),
    @(122, 6, 2024),
    @(123, 6, 2024),
    @(124, 7, 2024),
    @(125, 6, 2024),
    @(126, 6, 2024),
    @(127, 6, 2024),
    @(128, 6, 2024),
    @(129, 6, 2024),
    @(130, 5, 2024),
    @(131, 4, 2024),
    @(132, 5, 2024),
    @(133, 4, 2024),
    @(134, 5, 2024),
    @(135, 5, 2024),
    @(136, 4, 2024),
)

foreach ($item in $data) {
    $r = $item[0]
    $a = $item[1]
    $b = $item[2]
    $ws.Cells.Item($r, 1).Value2 = $a
    $ws.Cells.Item($r, 2).Value2 = $b
}

$win = $excel.ActiveWindow
$win.ScrollRow = 118
$win.ScrollColumn = 1

$ws.Range("D124").Select()
